$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New B:E and G values per row (F/"Win" column and A/date column are unchanged)
# Columns in each tuple: B(TB), C(d2S), D(K), E(IP), G(sum)
$newValues = @(
    @(3.230985683306322, 1.667794583268128, 0.8054896365839992, 8.660232485948974, 14.36450238910742),
    @(3.230985683306322, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 6.201049113329182),
    @(3.230985683306322, 1.667794583268128, 3.900430680208489, 0.496779210170732, 9.295990156953671),
    @(0.3048080303191223, 0.04240448674262143, 0.8054896365839992, 0.496779210170732, 1.649481363816475),
    @(0.127881588408715, 0.3127903958511391, 26.21740644021617, 8.660232485948974, 35.318310910425),
    @(1.459612070389937, 0.04240448674262143, 0.1575252929769615, 0.496779210170732, 2.156321060280252),
    @(0.3048080303191223, 0.3127903958511391, 3.900430680208489, 0.496779210170732, 5.014808316549482),
    @(3.230985683306322, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 6.201049113329182),
    @(0.127881588408715, 0.3127903958511391, 0.8054896365839992, 0.496779210170732, 1.742940831014585),
    @(0.3048080303191223, 10.29869402782916, 337.1190423067083, 8.660232485948974, 356.3827768508056),
    @(3.230985683306322, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 6.201049113329182),
    @(1.459612070389937, 1.667794583268128, 3.900430680208489, 0.496779210170732, 7.524616544037286),
    @(1.459612070389937, 1.667794583268128, 0.1575252929769615, 0.496779210170732, 3.781711156805759),
    @(1.459612070389937, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 4.429675500412797),
    @(1.459612070389937, 1.667794583268128, 3.900430680208489, 0.496779210170732, 7.524616544037286),
    @(3.230985683306322, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 6.201049113329182),
    @(0.3048080303191223, 0.3127903958511391, 0.1575252929769615, 0.496779210170732, 1.271902929317955),
    @(3.230985683306322, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 6.201049113329182),
    @(1.459612070389937, 1.667794583268128, 3.900430680208489, 0.496779210170732, 7.524616544037286),
    @(3.230985683306322, 1.667794583268128, 0.1575252929769615, 0.496779210170732, 5.553084769722144),
    @(3.230985683306322, 1.667794583268128, 0.1575252929769615, 0.496779210170732, 5.553084769722144),
    @(3.230985683306322, 114.8270160096505, 0.8054896365839992, 8.660232485948974, 127.5237238154898),
    @(3.230985683306322, 1.667794583268128, 0.8054896365839992, 8.660232485948974, 14.36450238910742)
)

$startRow = 2
for ($i = 0; $i -lt $newValues.Length; $i++) {
    $r = $startRow + $i
    $vals = $newValues[$i]
    $ws.Range("B$r").Value = $vals[0]
    $ws.Range("C$r").Value = $vals[1]
    $ws.Range("D$r").Value = $vals[2]
    $ws.Range("E$r").Value = $vals[3]
    $ws.Range("G$r").Value = $vals[4]
}

Write-Output "Updated $($newValues.Length) rows (B:E, G) for s_vals regen (filtered save games)."
